# Fruta / hortaliza, semanal
# Insert a new weekly record row for "Camote" (Zapallo) at row 59, shifting the
# existing rows 59:87 down to 60:88 (new report date 2021-09-13 / serial 44452).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 59 and below down by one to make room for the new record.
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new weekly observation.
$ws.Cells.Item(59, 1).Value  = 7
$ws.Cells.Item(59, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(59, 3).Value  = "Ñuble"
$ws.Cells.Item(59, 4).Value  = 44452
$ws.Cells.Item(59, 5).Value  = 16
$ws.Cells.Item(59, 6).Value  = 100112045
$ws.Cells.Item(59, 7).Value  = "Zapallo"
$ws.Cells.Item(59, 8).Value  = "Camote"
$ws.Cells.Item(59, 9).Value  = "1a (guarda)"
$ws.Cells.Item(59, 10).Value = 300
$ws.Cells.Item(59, 11).Value = 600
$ws.Cells.Item(59, 12).Value = 650
$ws.Cells.Item(59, 13).Value = 625
$ws.Cells.Item(59, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(59, 15).Value = "Región del Maule"
$ws.Cells.Item(59, 16).Value = 625
$ws.Cells.Item(59, 17).Value = 1
$ws.Cells.Item(59, 18).Value = "Hortaliza"
